$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# --- Title-case municipality/state name fixes ---
$ws.Range("B5").Value = 'Pabellón De Arteaga'
$ws.Range("B21").Value = 'Amatenango De La Frontera'
$ws.Range("B22").Value = 'Amatenango Del Valle'
$ws.Range("B30").Value = 'Comitán De Domínguez'
$ws.Range("B45").Value = 'Marqués De Comillas'
$ws.Range("B49").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B55").Value = 'San Cristóbal De Las Casas'
$ws.Range("B78").Value = 'Hidalgo Del Parral'
$ws.Range("B81").Value = 'San Francisco De Borja'
$ws.Range("A96").Value = 'Ciudad De México'
$ws.Range("A117").Value = 'Estado De México'
$ws.Range("B117").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B123").Value = 'Coacalco De Berriozábal'
$ws.Range("B126").Value = 'Ecatepec De Morelos'
$ws.Range("B129").Value = 'Ixtapan De La Sal'
$ws.Range("B133").Value = 'Naucalpan De Juárez'
$ws.Range("B138").Value = 'San Felipe Del Progreso'
$ws.Range("B144").Value = 'Tenango Del Valle'
$ws.Range("B146").Value = 'Tlalnepantla De Baz'
$ws.Range("B147").Value = 'Valle De Bravo'
$ws.Range("B148").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B149").Value = 'Villa Del Carbón'
$ws.Range("B156").Value = 'Apaseo El Alto'
$ws.Range("B157").Value = 'Apaseo El Grande'
$ws.Range("B165").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B175").Value = 'San Diego De La Unión'
$ws.Range("B178").Value = 'San Luis De La Paz'
$ws.Range("B181").Value = 'Valle De Santiago'
$ws.Range("B186").Value = 'Acapulco De Juárez'
$ws.Range("B188").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B191").Value = 'Atoyac De Álvarez'
$ws.Range("B193").Value = 'Buenavista De Cuéllar'
$ws.Range("B194").Value = 'Chilapa De Álvarez'
$ws.Range("B195").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B196").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B197").Value = 'Coyuca De Benítez'
$ws.Range("B198").Value = 'Coyuca De Catalán'
$ws.Range("B200").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B201").Value = 'Iguala De La Independencia'
$ws.Range("B202").Value = 'Zihuatanejo De Azueta'
$ws.Range("B203").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B209").Value = 'Taxco De Alarcón'
$ws.Range("B211").Value = 'Técpan De Galeana'
$ws.Range("B222").Value = 'Atotonilco El Grande'
$ws.Range("B225").Value = 'Huejutla De Reyes'
$ws.Range("B228").Value = 'Jacala De Ledezma'
$ws.Range("B232").Value = 'Mineral De La Reforma'
$ws.Range("B234").Value = 'Pachuca De Soto'
$ws.Range("B235").Value = 'Progreso De Obregón'
$ws.Range("B239").Value = 'Tenango De Doria'
$ws.Range("B240").Value = 'Tepehuacán De Guerrero'
$ws.Range("B241").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B244").Value = 'Tula De Allende'
$ws.Range("B249").Value = 'Atotonilco El Alto'
$ws.Range("B250").Value = 'Autlán De Navarro'
$ws.Range("B261").Value = 'San Juanito De Escobedo'
$ws.Range("B262").Value = 'Tamazula De Gordiano'
$ws.Range("B264").Value = 'Tepatitlán De Morelos'
$ws.Range("B271").Value = 'Zacoalco De Torres'
$ws.Range("B273").Value = 'Zapotlán El Grande'
$ws.Range("B313").Value = 'Puente De Ixtla'
$ws.Range("B325").Value = 'Montemorelos'
$ws.Range("B327").Value = 'San Nicolás De Los Garza'
$ws.Range("B329").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B330").Value = 'Coicoyán De Las Flores'
$ws.Range("B332").Value = 'Cuyamecalco Villa De Zaragoza'
$ws.Range("B333").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B334").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B335").Value = 'Ixtlán De Juárez'
$ws.Range("B336").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B339").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B340").Value = 'Oaxaca De Juárez'
$ws.Range("B355").Value = 'San Pedro El Alto'
$ws.Range("B375").Value = 'Teotitlán De Flores Magón'
$ws.Range("B376").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B387").Value = 'Ixcamilpa De Guerrero'
$ws.Range("B389").Value = 'Izúcar De Matamoros'
$ws.Range("B391").Value = 'Los Reyes De Juárez'
$ws.Range("B393").Value = 'Palmar De Bravo'
$ws.Range("B399").Value = 'Tetela De Ocampo'
$ws.Range("B401").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B409").Value = 'Amealco De Bonfil'
$ws.Range("B411").Value = 'Cadereyta De Montes'
$ws.Range("B415").Value = 'Jalpan De Serra'
$ws.Range("B417").Value = 'Pinal De Amoles'
$ws.Range("B426").Value = 'Armadillo De Los Infante'
$ws.Range("B427").Value = 'Axtla De Terrazas'
$ws.Range("B431").Value = 'Cerro De San Pedro'
$ws.Range("B432").Value = 'Ciudad Del Maíz'
$ws.Range("B442").Value = 'San Ciro De Acosta'
$ws.Range("B446").Value = 'Santa María Del Río'
$ws.Range("B451").Value = 'Tanquián De Escobedo'
$ws.Range("B453").Value = 'Villa De Arriaga'
$ws.Range("B454").Value = 'Villa De Ramos'
$ws.Range("B455").Value = 'Villa De Reyes'
$ws.Range("B495").Value = 'Soto La Marina'
$ws.Range("B503").Value = 'Acuamanala De Miguel Hidalgo'
$ws.Range("B507").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B515").Value = 'Amatlán De Los Reyes'
$ws.Range("B524").Value = 'Cosamaloapan De Carpio'
$ws.Range("B534").Value = 'Hueyapan De Ocampo'
$ws.Range("B542").Value = 'Juchique De Ferrer'
$ws.Range("B546").Value = 'Lerdo De Tejada'
$ws.Range("B549").Value = 'Martínez De La Torre'
$ws.Range("B558").Value = 'Poza Rica De Hidalgo'
$ws.Range("B564").Value = 'Sayula De Alemán'
$ws.Range("B565").Value = 'Soledad De Doblado'
$ws.Range("B567").Value = 'Tatahuicapan De Juárez'
$ws.Range("B596").Value = 'Villa De Cos'

# --- Numeric precision fix ---
$ws.Range("D71").Value = 0.09523809523809525

# --- Remove footer metadata rows 601-605 (and clean up dimension) ---
$ws.Range("A601:A605").EntireRow.Delete()
